{"js": "// Office.js (Word JavaScript API) implementation of the commit:\n//   \"Some code clean-up and commenting on the front-end code to help me\n//    understand it.\"\n//\n// Concretely this:\n//   1. Splits the \"for customers\" run so the (Word-managed) `_GoBack`\n//      bookmark sits right after \"...user types for c\" (i.e. where the\n//      author's cursor was left after their last edit), removing it from\n//      its old location at the very end of the document.\n//   2. Highlights the two \"Rename the customer table...\"/\"Add a new\n//      column...\" bullet points in yellow.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Move the `_GoBack` bookmark into the \"for customers\" sentence.\n// ---------------------------------------------------------------------\nconst splitResults = body.search(\"for c\", { matchCase: true, matchWholeWord: false });\nsplitResults.load(\"items\");\nawait context.sync();\n\nif (splitResults.items.length > 0) {\n  // Zero-width caret right after \"...for c\" (before \"ustomers...\").\n  const splitPoint = splitResults.items[0].getRange(\"End\");\n\n  // Word only ever keeps a single `_GoBack` bookmark (it marks the last\n  // edit position) -- drop the old one before inserting the new one so we\n  // don't end up with two.\n  context.document.deleteBookmark(\"_GoBack\");\n  splitPoint.insertBookmark(\"_GoBack\");\n}\n\n// ---------------------------------------------------------------------\n// 2) Yellow-highlight the two bullet points.\n// ---------------------------------------------------------------------\nconst highlightTexts = [\n  \"Rename the \\u2018customer\\u2019 table to \\u2018user\\u2019.\",\n  \"Add a new column in this database indicating user type. This will be an integer value. 1 will be used to identify employees, and 2 will be used to identify customers.\"\n];\n\nfor (const text of highlightTexts) {\n  const found = body.search(text, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    const para = found.items[0].paragraphs.getFirst();\n    para.font.highlightColor = \"Yellow\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) implementation of the commit:\n#   \"Some code clean-up and commenting on the front-end code to help me\n#    understand it.\"\n#\n# Concretely this:\n#   1. Splits the \"for customers\" run so the (Word-managed) `_GoBack`\n#      bookmark sits right after \"...user types for c\" (i.e. where the\n#      author's cursor was left after their last edit), moving it from its\n#      old location at the very end of the document.\n#   2. Highlights the two \"Rename the customer table...\"/\"Add a new\n#      column...\" bullet points in yellow.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Move the `_GoBack` bookmark into the \"for customers\" sentence.\n# ---------------------------------------------------------------------\n$splitRange = $d.Content\n$found = $splitRange.Find.Execute(\"for c\")\nif ($found) {\n  $splitRange.Collapse(0)   # wdCollapseEnd -- zero-width caret after \"for c\"\n  # Bookmark names are unique -- re-adding \"_GoBack\" here both creates it\n  # at the new spot and removes it from wherever it used to be.\n  $d.Bookmarks.Add(\"_GoBack\", $splitRange) | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 2) Yellow-highlight the two bullet points.\n# ---------------------------------------------------------------------\n$wdYellow = 7\n\nforeach ($p in $d.Paragraphs) {\n  $t = $p.Range.Text\n  if ($t.StartsWith(\"Rename the\") -or $t.StartsWith(\"Add a new column in this database indicating user type\")) {\n    $p.Range.Font.HighlightColorIndex = $wdYellow\n  }\n}\n"}
